$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 492.44446
$ws.Range("I8").Value = 492.44446
$ws.Range("K8").Value = 1477.33338
$ws.Range("M8").Value = -1338.33338
$ws.Range("H33").Value = 19518.871
$ws.Range("I33").Value = 22368.73
$ws.Range("J33").Value = 4699.6
$ws.Range("K33").Value = 22368.73
$ws.Range("L33").Value = 4699.6
$ws.Range("M33").Value = -22139.73
$ws.Range("N33").Value = -5157.6
$ws.Range("H40").Value = 1849.7142
$ws.Range("I40").Value = 1933
$ws.Range("J40").Value = 1787.25
$ws.Range("K40").Value = 1933
$ws.Range("L40").Value = 1787.25
$ws.Range("M40").Value = -1758
$ws.Range("N40").Value = -2137.25
$ws.Range("H86").Value = 76926344
$ws.Range("I86").Value = 76926344
$ws.Range("K86").Value = 76926344
$ws.Range("M86").Value = -76925221
$ws.Range("H89").Value = 76926344
$ws.Range("I89").Value = 76926344
$ws.Range("K89").Value = 384631720
$ws.Range("M89").Value = -384626104
$ws.Range("H98").Value = 3574.125
$ws.Range("I98").Value = 2720.1
$ws.Range("J98").Value = 4997.5
$ws.Range("K98").Value = 2720.1
$ws.Range("L98").Value = 4997.5
$ws.Range("M98").Value = -1222.1
$ws.Range("N98").Value = -7993.5
$ws.Range("H106").Value = 4221.1665
$ws.Range("I106").Value = 2786.7273
$ws.Range("K106").Value = 2786.7273
$ws.Range("M106").Value = -2155.7273
$ws.Range("H122").Value = 3574.125
$ws.Range("I122").Value = 2720.1
$ws.Range("J122").Value = 4997.5
$ws.Range("K122").Value = 8160.299999999999
$ws.Range("L122").Value = 14992.5
$ws.Range("M122").Value = -5710.299999999999
$ws.Range("N122").Value = -19892.5
$ws.Range("H132").Value = 1388.7709
$ws.Range("I132").Value = 1335.8292
$ws.Range("K132").Value = 4007.487599999999
$ws.Range("M132").Value = -1477.487599999999
$ws.Range("H135").Value = 1848.6316
$ws.Range("I135").Value = 1432.2727
$ws.Range("J135").Value = 2421.125
$ws.Range("K135").Value = 12890.4543
$ws.Range("L135").Value = 21790.125
$ws.Range("M135").Value = -10355.4543
$ws.Range("N135").Value = -26860.125

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 8166.1665
$ws.Range("I86").Value = 1499.75
$ws.Range("J86").Value = 11499.375
$ws.Range("K86").Value = 1499.75
$ws.Range("L86").Value = 11499.375
$ws.Range("M86").Value = -376.75
$ws.Range("N86").Value = -13745.375
$ws.Range("H89").Value = 8166.1665
$ws.Range("I89").Value = 1499.75
$ws.Range("J89").Value = 11499.375
$ws.Range("K89").Value = 7498.75
$ws.Range("L89").Value = 57496.875
$ws.Range("M89").Value = -1882.75
$ws.Range("N89").Value = -68728.875
$ws.Range("H94").Value = 2012.7273
$ws.Range("I94").Value = 1345.1177
$ws.Range("J94").Value = 4282.6
$ws.Range("K94").Value = 1345.1177
$ws.Range("L94").Value = 4282.6
$ws.Range("M94").Value = -894.1177
$ws.Range("N94").Value = -5184.6
$ws.Range("H105").Value = 19230.54
$ws.Range("I105").Value = 21385.2
$ws.Range("J105").Value = 17883.875
$ws.Range("K105").Value = 21385.2
$ws.Range("L105").Value = 17883.875
$ws.Range("M105").Value = -19638.2
$ws.Range("N105").Value = -21377.875
$ws.Range("H134").Value = 32145570
$ws.Range("I134").Value = 2231.625
$ws.Range("K134").Value = 6694.875
$ws.Range("M134").Value = -4159.875

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 35717660
$ws.Range("I16").Value = 47620212
$ws.Range("K16").Value = 47620212
$ws.Range("M16").Value = -47619925
$ws.Range("H31").Value = 2935.5715
$ws.Range("I31").Value = 3651.7334
$ws.Range("K31").Value = 3651.7334
$ws.Range("M31").Value = -3356.7334
$ws.Range("H34").Value = 2935.5715
$ws.Range("I34").Value = 3651.7334
$ws.Range("K34").Value = 3651.7334
$ws.Range("M34").Value = -3449.7334
$ws.Range("H113").Value = 35717660
$ws.Range("I113").Value = 47620212
$ws.Range("K113").Value = 47620212
$ws.Range("M113").Value = -47618042
$ws.Range("H122").Value = 2808.36
$ws.Range("I122").Value = 2708.6365
$ws.Range("K122").Value = 8125.9095
$ws.Range("M122").Value = -5675.9095

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 6389413.5
$ws.Range("J9").Value = 9287433
$ws.Range("L9").Value = 27862299
$ws.Range("N9").Value = -27862747
$ws.Range("H31").Value = 1750
$ws.Range("I31").Value = 1750
$ws.Range("K31").Value = 5250
$ws.Range("M31").Value = -4962
$ws.Range("H56").Value = 12825062
$ws.Range("I56").Value = 12825062
$ws.Range("K56").Value = 12825062
$ws.Range("M56").Value = -12824532
$ws.Range("H59").Value = 20000
$ws.Range("I59").Value = 20000
$ws.Range("J59").Value = 0
$ws.Range("K59").Value = 60000
$ws.Range("L59").Value = 0
$ws.Range("M59").Value = -59460
$ws.Range("N59").ClearContents()
$ws.Range("H61").Value = 7081.6665
$ws.Range("J61").Value = 10122.5
$ws.Range("L61").Value = 30367.5
$ws.Range("N61").Value = -30797.5
$ws.Range("H94").Value = 16996.834
$ws.Range("I94").Value = 10000
$ws.Range("K94").Value = 30000
$ws.Range("M94").Value = -29324
$ws.Range("H105").Value = 13383.167
$ws.Range("J105").Value = 14059.8
$ws.Range("L105").Value = 42179.39999999999
$ws.Range("N105").Value = -47421.39999999999
$ws.Range("H119").Value = 111126550
$ws.Range("I119").Value = 166677330
$ws.Range("J119").Value = 24995
$ws.Range("K119").Value = 500031990
$ws.Range("L119").Value = 74985
$ws.Range("M119").Value = -500027152
$ws.Range("N119").Value = -84661
$ws.Range("H132").Value = 2089.4211
$ws.Range("I132").Value = 1969.1538
$ws.Range("J132").Value = 2350
$ws.Range("K132").Value = 17722.3842
$ws.Range("L132").Value = 21150
$ws.Range("M132").Value = -15192.3842
$ws.Range("N132").Value = -26210

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 62501944
$ws.Range("I102").Value = 83335190
$ws.Range("K102").Value = 83335190
$ws.Range("M102").Value = -83333568
$ws.Range("H122").Value = 1990.5834
$ws.Range("I122").Value = 1720.7778
$ws.Range("J122").Value = 2800
$ws.Range("K122").Value = 5162.3334
$ws.Range("L122").Value = 8400
$ws.Range("M122").Value = -2712.3334
$ws.Range("N122").Value = -13300
$ws.Range("H126").Value = 2261.6667
$ws.Range("I126").Value = 2166
$ws.Range("K126").Value = 6498
$ws.Range("M126").Value = -4028
$ws.Range("H132").Value = 755959.3
$ws.Range("J132").Value = 988143.1
$ws.Range("L132").Value = 2964429.3
$ws.Range("N132").Value = -2969489.3

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 99999
$ws.Range("J2").Value = 99999
$ws.Range("L2").Value = 99999
$ws.Range("N2").Value = -100223
$ws.Range("H7").Value = 7006.5386
$ws.Range("I7").Value = 3409.1
$ws.Range("K7").Value = 3409.1
$ws.Range("M7").Value = -3297.1
$ws.Range("H40").Value = 2453.4783
$ws.Range("I40").Value = 2453.4783
$ws.Range("K40").Value = 2453.4783
$ws.Range("M40").Value = -2317.4783
$ws.Range("H46").Value = 9843.277
$ws.Range("I46").Value = 16197.429
$ws.Range("J46").Value = 5799.727
$ws.Range("K46").Value = 16197.429
$ws.Range("L46").Value = 5799.727
$ws.Range("M46").Value = -16009.429
$ws.Range("N46").Value = -6175.727
$ws.Range("H94").Value = 39998.5
$ws.Range("J94").Value = 39998.5
$ws.Range("L94").Value = 39998.5
$ws.Range("N94").Value = -41350.5
$ws.Range("H126").Value = 7006.5386
$ws.Range("I126").Value = 3409.1
$ws.Range("K126").Value = 10227.3
$ws.Range("M126").Value = -7757.299999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 2607.3333
$ws.Range("I96").Value = 3028.8
$ws.Range("K96").Value = 3028.8
$ws.Range("M96").Value = -1655.8
$ws.Range("H110").Value = 50644
$ws.Range("J110").Value = 50644
$ws.Range("L110").Value = 50644
$ws.Range("N110").Value = -58824
$ws.Range("H113").Value = 870.4
$ws.Range("I113").Value = 126
$ws.Range("K113").Value = 378
$ws.Range("M113").Value = 1792
$ws.Range("H122").Value = 1392.8334
$ws.Range("I122").Value = 1499.1428
$ws.Range("J122").Value = 1020.75
$ws.Range("K122").Value = 4497.428400000001
$ws.Range("L122").Value = 3062.25
$ws.Range("M122").Value = -2047.428400000001
$ws.Range("N122").Value = -7962.25
$ws.Range("H132").Value = 2819.7334
$ws.Range("I132").Value = 2391.8572
$ws.Range("K132").Value = 7175.571599999999
$ws.Range("M132").Value = -4645.571599999999
